$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells that are being updated, to preserve exact
# string formatting (e.g. trailing zeros, multi-dot thousand separators) instead of
# having Excel auto-convert them to numbers.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.949.04"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "1.891.27"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("D4").Value = "1.020"
$ws.Range("E4").Value = "  +1.76%  "
$ws.Range("D5").Value = "336.05"
$ws.Range("E5").Value = "  +1.48%  "
$ws.Range("D6").Value = "1.017"
$ws.Range("E6").Value = "  +1.56%  "
$ws.Range("D7").Value = "0.4690"
$ws.Range("E7").Value = "  -1.21%  "
$ws.Range("D8").Value = "0.3919"
$ws.Range("E8").Value = "  -1.68%  "
$ws.Range("D9").Value = "47.49"
$ws.Range("E9").Value = "  -0.95%  "
$ws.Range("D10").Value = "0.08029"
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").Value = "1.019"
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").Value = "21.77"
$ws.Range("E12").Value = "  -0.22%  "
$ws.Range("D13").Value = "1.890.10"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").Value = "5.960"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "7.107"
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").Value = "1.020"
$ws.Range("E16").Value = "  +1.84%  "
$ws.Range("D17").Value = "0.06799"
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.00001051"
$ws.Range("E18").Value = "  +0.96%  "
$ws.Range("B19").Value = "Litecoin"
$ws.Range("C19").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D19").Value = "87.36"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "17.15"
$ws.Range("E20").Value = "  -0.89%  "
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("D22").Value = "27.968.27"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("D23").Value = "5.508"
$ws.Range("E23").Value = "  -0.22%  "
$ws.Range("D24").Value = "10.97"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").Value = "2.347"
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("D26").Value = "2.120.29"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "159.30"
$ws.Range("E27").Value = "  +2.49%  "
$ws.Range("D28").Value = "20.05"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").Value = "2.076"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("D30").Value = "5.430"
$ws.Range("E30").Value = "  -2.82%  "
$ws.Range("D31").Value = "121.70"
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("D32").Value = "0.9682"
$ws.Range("E32").Value = "  -0.22%  "
$ws.Range("D33").Value = "0.09509"
$ws.Range("E33").Value = "  -0.47%  "
$ws.Range("D34").Value = "3.672"
$ws.Range("E34").Value = "  +1.32%  "
$ws.Range("D35").Value = "1.386"
$ws.Range("E35").Value = "  -5.56%  "
$ws.Range("D36").Value = "5.345"
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("D37").Value = "0.06120"
$ws.Range("E37").Value = "  -0.04%  "
$ws.Range("D38").Value = "0.02246"
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("D39").Value = "1.218"
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("D40").Value = "8.092"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").Value = "0.5979"
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("D42").Value = "0.1890"
$ws.Range("E42").Value = "  -0.44%  "
$ws.Range("D43").Value = "10.33"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  +1.36%  "
$ws.Range("D45").Value = "0.5689"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("D46").Value = "12.20"
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "3.404"
$ws.Range("E47").Value = "  -0.12%  "
$ws.Range("D48").Value = "1.935"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.06932"
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "113.89"
$ws.Range("E50").Value = "  +2.80%  "
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "1.070"
$ws.Range("E51").Value = "  +0.12%  "

# Restore default (no explicit number-format) style on the Price cells we touched,
# matching the original workbook styling (General format, no style index).
foreach ($c in $priceCells) {
    $ws.Range($c).Style = "Normal"
}
